$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Notes": collapse the multi-range selection down to A4 only
# ------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Activate() | Out-Null
$wsNotes.Range("A4").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "CARS": collapse the multi-range selection down to A43 only
# ------------------------------------------------------------------
$wsCars = $wb.Worksheets.Item("CARS")
$wsCars.Activate() | Out-Null
$wsCars.Range("A43").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "LGV": collapse the multi-range selection down to A43 only
# ------------------------------------------------------------------
$wsLgv = $wb.Worksheets.Item("LGV")
$wsLgv.Activate() | Out-Null
$wsLgv.Range("A43").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "HGV": add a totals column (I), hide the yearly C:H columns,
# tighten row heights and move the view/selection to I1
# ------------------------------------------------------------------
$wsHgv = $wb.Worksheets.Item("HGV")
$wsHgv.Activate() | Out-Null

# Add SUM formulas in column I for each data row (4-43); Excel will
# auto-adjust the relative references row by row.
$wsHgv.Range("I4:I43").Formula = "=SUM(C4:H4)"

# Match the formatting (style) of column A for the new column I.
$wsHgv.Range("A4:A43").Copy() | Out-Null
$wsHgv.Range("I4:I43").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Tighten the row heights for the data rows.
$wsHgv.Range("A4:A43").EntireRow.RowHeight = 13.8

# Hide the now-redundant per-year columns C:H.
$wsHgv.Range("C:H").EntireColumn.Hidden = $true

# Move the view back to the top and select the new total cell.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsHgv.Range("I1").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "OTHER VEHICLES": collapse the multi-range selection down to
# D22 only
# ------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("OTHER VEHICLES")
$wsOther.Activate() | Out-Null
$wsOther.Range("D22").Select() | Out-Null

# Restore HGV as the active/selected sheet, matching the workbook's
# original active tab.
$wsHgv.Activate() | Out-Null

Write-Host "done"
